$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2895.7778
$ws.Range("I40").Value = 2572.6
$ws.Range("J40").Value = 3299.75
$ws.Range("K40").Value = 2572.6
$ws.Range("L40").Value = 3299.75
$ws.Range("M40").Value = -2397.6
$ws.Range("N40").Value = -3649.75
$ws.Range("H64").Value = 8147.0605
$ws.Range("I64").Value = 4317
$ws.Range("K64").Value = 4317
$ws.Range("M64").Value = -4069
$ws.Range("H67").Value = 8147.0605
$ws.Range("I67").Value = 4317
$ws.Range("K67").Value = 4317
$ws.Range("M67").Value = -3459
$ws.Range("H92").Value = 1285.08
$ws.Range("I92").Value = 1276.8889
$ws.Range("K92").Value = 1276.8889
$ws.Range("M92").Value = -28.88889999999992
$ws.Range("H94").Value = 4512
$ws.Range("I94").Value = 3449.4
$ws.Range("K94").Value = 3449.4
$ws.Range("M94").Value = -2998.4
$ws.Range("H106").Value = 5111.3716
$ws.Range("I106").Value = 5954.9473
$ws.Range("J106").Value = 4109.625
$ws.Range("K106").Value = 5954.9473
$ws.Range("L106").Value = 4109.625
$ws.Range("M106").Value = -5323.9473
$ws.Range("N106").Value = -5371.625
$ws.Range("H112").Value = 2512.2307
$ws.Range("J112").Value = 4654.75
$ws.Range("L112").Value = 13964.25
$ws.Range("N112").Value = -16180.25
$ws.Range("H116").Value = 3828.4285
$ws.Range("J116").Value = 3759.8
$ws.Range("L116").Value = 3759.8
$ws.Range("N116").Value = -10643.8
$ws.Range("H137").Value = 1651.375
$ws.Range("I137").Value = 1157.2222
$ws.Range("K137").Value = 3471.6666
$ws.Range("M137").Value = -921.6665999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 7357
$ws.Range("J63").Value = 9166.666999999999
$ws.Range("L63").Value = 9166.666999999999
$ws.Range("N63").Value = -10538.667
$ws.Range("H66").Value = 7357
$ws.Range("J66").Value = 9166.666999999999
$ws.Range("L66").Value = 45833.335
$ws.Range("N66").Value = -52697.335
$ws.Range("H102").Value = 4858.8667
$ws.Range("I102").Value = 2688.6
$ws.Range("K102").Value = 2688.6
$ws.Range("M102").Value = -1066.6
$ws.Range("H132").Value = 1732
$ws.Range("I132").Value = 1766.8572
$ws.Range("K132").Value = 5300.571599999999
$ws.Range("M132").Value = -2770.571599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2683.9583
$ws.Range("I86").Value = 2089.1765
$ws.Range("J86").Value = 4128.4287
$ws.Range("K86").Value = 2089.1765
$ws.Range("L86").Value = 4128.4287
$ws.Range("M86").Value = -966.1765
$ws.Range("N86").Value = -6374.4287
$ws.Range("H89").Value = 2683.9583
$ws.Range("I89").Value = 2089.1765
$ws.Range("J89").Value = 4128.4287
$ws.Range("K89").Value = 10445.8825
$ws.Range("L89").Value = 20642.1435
$ws.Range("M89").Value = -4829.8825
$ws.Range("N89").Value = -31874.1435
$ws.Range("H134").Value = 4232.8184
$ws.Range("I134").Value = 5260.1665
$ws.Range("K134").Value = 15780.4995
$ws.Range("M134").Value = -13245.4995
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3178.0789
$ws.Range("I31").Value = 1656.1666
$ws.Range("J31").Value = 3880.5
$ws.Range("K31").Value = 1656.1666
$ws.Range("L31").Value = 3880.5
$ws.Range("M31").Value = -1361.1666
$ws.Range("N31").Value = -4470.5
$ws.Range("H34").Value = 3178.0789
$ws.Range("I34").Value = 1656.1666
$ws.Range("J34").Value = 3880.5
$ws.Range("K34").Value = 1656.1666
$ws.Range("L34").Value = 3880.5
$ws.Range("M34").Value = -1454.1666
$ws.Range("N34").Value = -4284.5
$ws.Range("H132").Value = 1552.5
$ws.Range("I132").Value = 1467
$ws.Range("K132").Value = 4401
$ws.Range("M132").Value = -1871
$ws.Range("H141").Value = 209999.67
$ws.Range("J141").Value = 209999.67
$ws.Range("L141").Value = 209999.67
$ws.Range("N141").Value = -220359.67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 105333.336
$ws.Range("I56").Value = 105333.336
$ws.Range("K56").Value = 105333.336
$ws.Range("M56").Value = -104803.336
$ws.Range("H131").Value = 1698.3846
$ws.Range("J131").Value = 2061.3572
$ws.Range("L131").Value = 6184.071599999999
$ws.Range("N131").Value = -16264.0716
$ws.Range("H132").Value = 1956.2632
$ws.Range("I132").Value = 1190
$ws.Range("K132").Value = 10710
$ws.Range("M132").Value = -8180
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4425.3076
$ws.Range("I80").Value = 2563.7778
$ws.Range("J80").Value = 5410.8237
$ws.Range("K80").Value = 2563.7778
$ws.Range("L80").Value = 5410.8237
$ws.Range("M80").Value = -1565.7778
$ws.Range("N80").Value = -7406.8237
$ws.Range("H83").Value = 4425.3076
$ws.Range("I83").Value = 2563.7778
$ws.Range("J83").Value = 5410.8237
$ws.Range("K83").Value = 12818.889
$ws.Range("L83").Value = 27054.1185
$ws.Range("M83").Value = -7826.888999999999
$ws.Range("N83").Value = -37038.1185
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1823.5834
$ws.Range("J46").Value = 1973.4445
$ws.Range("L46").Value = 1973.4445
$ws.Range("N46").Value = -2349.4445
$ws.Range("H68").Value = 8166.3335
$ws.Range("I68").Value = 5749
$ws.Range("J68").Value = 9375
$ws.Range("K68").Value = 5749
$ws.Range("L68").Value = 9375
$ws.Range("M68").Value = -5000
$ws.Range("N68").Value = -10873
$ws.Range("H71").Value = 8166.3335
$ws.Range("I71").Value = 5749
$ws.Range("J71").Value = 9375
$ws.Range("K71").Value = 28745
$ws.Range("L71").Value = 46875
$ws.Range("M71").Value = -25001
$ws.Range("N71").Value = -54363
$ws.Range("H93").Value = 3101.3057
$ws.Range("I93").Value = 575.63635
$ws.Range("K93").Value = 575.63635
$ws.Range("M93").Value = 672.36365
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25500
$ws.Range("H36").Value = 25000
$ws.Range("J36").Value = 25000
$ws.Range("L36").Value = 25000
$ws.Range("N36").Value = -25500
$ws.Range("H96").Value = 4444
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H126").Value = 2756.35
$ws.Range("I126").Value = 2291.6875
$ws.Range("J126").Value = 4615
$ws.Range("K126").Value = 6875.0625
$ws.Range("L126").Value = 13845
$ws.Range("M126").Value = -4405.0625
$ws.Range("N126").Value = -18785
